# Edit the class-diagram PowerPoint:
#  - rename a handful of diagram shape labels on slide 1 (AddressBook -> Giatros
#    domain rename: Person -> Patient, Tag -> Allergy, etc.)
#  - refresh the cached "datetimeFigureOut" footer date (12/5/2018 -> 4/15/2019)
#    everywhere it is stamped: the slide master, every slide layout, and the
#    notes master.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Diagram label renames on slide 1
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if (-not $sh.HasTextFrame) { continue }
    if (-not $sh.TextFrame.HasText) { continue }

    $t = $sh.TextFrame.TextRange.Text

    if ($t -eq "VersionedAddressBook") {
        $sh.TextFrame.TextRange.Text = "VersionedGiatrosBook"
    }
    elseif ($t -eq "UniquePersonList") {
        $sh.TextFrame.TextRange.Text = "UniquePatientList"
    }
    elseif ($t -eq "Person") {
        $sh.TextFrame.TextRange.Text = "Patient"
    }
    elseif ($t -eq "UniqueTagList") {
        $sh.TextFrame.TextRange.Text = "UniqueAllergyList"
    }
    elseif ($t -eq "Tag") {
        $sh.TextFrame.TextRange.Text = "Allergy"
    }
}

# ---------------------------------------------------------------------------
# 2. Refresh the cached date-footer field text (12/5/2018 -> 4/15/2019)
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if (-not $sh.HasTextFrame) { continue }
        if (-not $sh.TextFrame.HasText) { continue }
        if ($sh.TextFrame.TextRange.Text -eq "12/5/2018") {
            $sh.TextFrame.TextRange.Text = "4/15/2019"
        }
    }
}

# Slide master
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout hanging off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShapes $layouts.Item($j).Shapes
}

# Notes master
Update-DateShapes $p.NotesMaster.Shapes
